# Update "想去人数" (F column) counts for rows 2-10 on both the "展览"
# sheet and the "全部类型" sheet (they carry duplicate data).

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 2948
    3  = 732
    4  = 103
    5  = 6716
    6  = 1686
    7  = 20
    8  = 31
    9  = 57
    10 = 114
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
